# Update NATMI ligand-receptor TPM stats (Tgfb2-Tgfbr1) on Sheet1 with
# newly recomputed values (rows 2-10, columns E-T as applicable).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8606349999999999
$ws.Range("H2").Value = 2.581905
$ws.Range("I2").Value = 0.0262626340301864
$ws.Range("J2").Value = 0.0262626340301864
$ws.Range("M2").Value = 5.482938999999999
$ws.Range("N2").Value = 16.448817
$ws.Range("O2").Value = 0.1472261722051079
$ws.Range("P2").Value = 0.147226172205108
$ws.Range("Q2").Value = 4.718809206264999
$ws.Range("R2").Value = 42.46928285638499
$ws.Range("S2").Value = 0.003866547080287951
$ws.Range("T2").Value = 0.003866547080287951
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8606349999999999
$ws.Range("H3").Value = 2.581905
$ws.Range("I3").Value = 0.0262626340301864
$ws.Range("J3").Value = 0.0262626340301864
$ws.Range("O3").Value = 0.5993885906243068
$ws.Range("P3").Value = 0.5993885906243068
$ws.Range("Q3").Value = 19.21126085943333
$ws.Range("R3").Value = 172.9013477349
$ws.Range("S3").Value = 0.01574152319743538
$ws.Range("T3").Value = 0.01574152319743538
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8606349999999999
$ws.Range("H4").Value = 2.581905
$ws.Range("I4").Value = 0.0262626340301864
$ws.Range("J4").Value = 0.0262626340301864
$ws.Range("M4").Value = 9.436472999999999
$ws.Range("N4").Value = 28.309419
$ws.Range("O4").Value = 0.2533852371705853
$ws.Range("P4").Value = 0.2533852371705853
$ws.Range("Q4").Value = 8.121358940354998
$ws.Range("R4").Value = 73.09223046319499
$ws.Range("S4").Value = 0.006654563752463065
$ws.Range("T4").Value = 0.006654563752463065
$ws.Range("H5").Value = 58.40949000000001
$ws.Range("I5").Value = 0.5941299388474139
$ws.Range("J5").Value = 0.5941299388474139
$ws.Range("M5").Value = 5.482938999999999
$ws.Range("N5").Value = 16.448817
$ws.Range("O5").Value = 0.1472261722051079
$ws.Range("P5").Value = 0.147226172205108
$ws.Range("Q5").Value = 106.75189023037
$ws.Range("R5").Value = 960.76701207333
$ws.Range("S5").Value = 0.08747147668895962
$ws.Range("T5").Value = 0.08747147668895963
$ws.Range("H6").Value = 58.40949000000001
$ws.Range("I6").Value = 0.5941299388474139
$ws.Range("J6").Value = 0.5941299388474139
$ws.Range("O6").Value = 0.5993885906243068
$ws.Range("P6").Value = 0.5993885906243068
$ws.Range("R6").Value = 3911.4837848442
$ws.Range("S6").Value = 0.356114706693457
$ws.Range("T6").Value = 0.356114706693457
$ws.Range("H7").Value = 58.40949000000001
$ws.Range("I7").Value = 0.5941299388474139
$ws.Range("J7").Value = 0.5941299388474139
$ws.Range("M7").Value = 9.436472999999999
$ws.Range("N7").Value = 28.309419
$ws.Range("O7").Value = 0.2533852371705853
$ws.Range("P7").Value = 0.2533852371705853
$ws.Range("Q7").Value = 183.72652510959
$ws.Range("R7").Value = 1653.53872598631
$ws.Range("S7").Value = 0.1505437554649973
$ws.Range("T7").Value = 0.1505437554649973
$ws.Range("G8").Value = 12.439858
$ws.Range("H8").Value = 37.319574
$ws.Range("I8").Value = 0.3796074271223998
$ws.Range("J8").Value = 0.3796074271223997
$ws.Range("M8").Value = 5.482938999999999
$ws.Range("N8").Value = 16.448817
$ws.Range("O8").Value = 0.1472261722051079
$ws.Range("P8").Value = 0.147226172205108
$ws.Range("Q8").Value = 68.20698258266198
$ws.Range("R8").Value = 613.8628432439579
$ws.Range("S8").Value = 0.05588814843586039
$ws.Range("T8").Value = 0.05588814843586039
$ws.Range("G9").Value = 12.439858
$ws.Range("H9").Value = 37.319574
$ws.Range("I9").Value = 0.3796074271223998
$ws.Range("J9").Value = 0.3796074271223997
$ws.Range("O9").Value = 0.5993885906243068
$ws.Range("P9").Value = 0.5993885906243068
$ws.Range("Q9").Value = 277.6849153152133
$ws.Range("R9").Value = 2499.164237836919
$ws.Range("S9").Value = 0.2275323607334145
$ws.Range("T9").Value = 0.2275323607334144
$ws.Range("G10").Value = 12.439858
$ws.Range("H10").Value = 37.319574
$ws.Range("I10").Value = 0.3796074271223998
$ws.Range("J10").Value = 0.3796074271223997
$ws.Range("M10").Value = 9.436472999999999
$ws.Range("N10").Value = 28.309419
$ws.Range("O10").Value = 0.2533852371705853
$ws.Range("P10").Value = 0.2533852371705853
$ws.Range("Q10").Value = 117.388384140834
$ws.Range("R10").Value = 1056.495457267506
$ws.Range("S10").Value = 0.09618691795312494
$ws.Range("T10").Value = 0.09618691795312492
